$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Improving Firefox Stability in the Enterprise by Reducing DLL Injection"
$ws.Range("A3").Value = "Launching Interop 2025"
$ws.Range("A4").Value = "Introducing Uniffi for React Native: Rust-Powered Turbo Modules"
$ws.Range("A5").Value = "Llamafile v0.8.14: a new UI, performance gains, and more"
$ws.Range("A6").Value = "0Din: A GenAI Bug Bounty Program – Securing Tomorrow’s AI Together"
$ws.Range("A7").Value = "Announcing Official Puppeteer Support for Firefox"
$ws.Range("A8").Value = "Snapshots for IPC Fuzzing"
$ws.Range("A9").Value = "Sponsoring sqlite-vec to enable more powerful Local AI applications"
$ws.Range("A10").Value = "Experimenting with local alt text generation in Firefox Nightly"
$ws.Range("A11").Value = "Llamafile’s progress, four months in"
$ws.Range("A12").Value = "Porting a cross-platform GUI application to Rust"
